$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.092.95"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").Value = "1.858.74"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'232.84"
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = "  -3.20%  "
$ws.Range("D8").Value = "'0.2811"
$ws.Range("E8").Value = "  -3.70%  "
$ws.Range("D9").Value = "'0.06529"
$ws.Range("E9").Value = "  -4.43%  "
$ws.Range("D10").Value = "'19.50"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "'0.07808"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'96.25"
$ws.Range("E12").Value = "  -8.08%  "
$ws.Range("D13").Value = "1.856.02"
$ws.Range("E13").Value = "  -5.26%  "
$ws.Range("D14").Value = "'5.117"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "'0.6630"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "'279.81"
$ws.Range("E16").Value = "  -5.48%  "
$ws.Range("D17").Value = "30.121.08"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'5.460"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "'12.55"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").Value = "2.095.88"
$ws.Range("E21").Value = "  -4.95%  "
$ws.Range("D22").Value = "'0.000007204"
$ws.Range("E22").Value = "  -5.76%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'6.104"
$ws.Range("E24").Value = "  -5.84%  "
$ws.Range("D25").Value = "'9.296"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("D26").Value = "'165.53"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").Value = "'18.80"
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("D28").Value = "'1.906"
$ws.Range("E28").Value = "  -10.65%  "
$ws.Range("D29").Value = "'1.332"
$ws.Range("E29").Value = "  -5.09%  "
$ws.Range("D30").Value = "'0.09525"
$ws.Range("E30").Value = "  -6.52%  "
$ws.Range("D31").Value = "'4.395"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D32").Value = "'1.463"
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "'4.084"
$ws.Range("E33").Value = "  -6.67%  "
$ws.Range("D34").Value = "'0.04643"
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("D35").Value = "'0.7002"
$ws.Range("E35").Value = "  -6.10%  "
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'0.01849"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").Value = "'6.259"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").Value = "'2.510"
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("D41").Value = "'72.55"
$ws.Range("E41").Value = "  -5.84%  "
$ws.Range("D42").Value = "'0.8520"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "'1.913"
$ws.Range("E43").Value = "  -7.26%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'103.64"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "'0.4136"
$ws.Range("E46").Value = "  -5.91%  "
$ws.Range("D47").Value = "'990.44"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "'7.156"
$ws.Range("E48").Value = "  -6.04%  "
$ws.Range("D49").Value = "'9.180"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").Value = "'0.1134"
$ws.Range("E51").Value = "  -6.87%  "
